$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

$ws1.Range("B2").Value = "Wed Dec 20 12:37:16 EST 2023"
$ws1.Range("B3").Value = "Wed Dec 20 12:37:29 EST 2023"
$ws1.Range("B4").Value = "Wed Dec 20 12:37:41 EST 2023"
$ws1.Range("B5").Value = "Wed Dec 20 12:37:52 EST 2023"
$ws1.Range("B6").Value = "Wed Dec 20 12:38:05 EST 2023"
$ws1.Range("B7").Value = "Wed Dec 20 12:38:18 EST 2023"
$ws1.Range("B8").Value = "Wed Dec 20 12:38:29 EST 2023"
$ws1.Range("B9").Value = "Wed Dec 20 12:38:40 EST 2023"
$ws1.Range("B10").Value = "Wed Dec 20 12:38:51 EST 2023"
$ws1.Range("B13").Value = "Wed Dec 20 12:39:02 EST 2023"

$ws2.Range("B2").Value = "Wed Dec 20 12:39:14 EST 2023"
$ws2.Range("B3").Value = "Wed Dec 20 12:39:25 EST 2023"
$ws2.Range("B4").Value = "Wed Dec 20 12:39:36 EST 2023"
$ws2.Range("B5").Value = "Wed Dec 20 12:39:47 EST 2023"
$ws2.Range("B6").Value = "Wed Dec 20 12:39:58 EST 2023"
$ws2.Range("B7").Value = "Wed Dec 20 12:40:09 EST 2023"
$ws2.Range("B8").Value = "Wed Dec 20 12:40:20 EST 2023"
$ws2.Range("B9").Value = "Wed Dec 20 12:40:31 EST 2023"
$ws2.Range("B14").Value = "Wed Dec 20 12:40:42 EST 2023"
$ws2.Range("B15").Value = "Wed Dec 20 12:40:53 EST 2023"
$ws2.Range("B16").Value = "Wed Dec 20 12:41:04 EST 2023"
$ws2.Range("B17").Value = "Wed Dec 20 12:41:15 EST 2023"
